# Update the "Förändrad" (column C) date value for rows 2-9 from
# 2023-11-03 (45233) to 2023-11-13 (45243), as in the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 9; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45233) {
        $cell.Value2 = 45243
    }
}
